# Update the "450 DSA Questions" progress tracker on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 19-24 (Array section): mark progress as "yes"
$ws.Range("C19:C24").Value = "yes"

# Rows 42-43 (Matrix section): clear the placeholder marker entirely
$ws.Range("C42:C43").ClearContents()

# Row 139 (LinkedList section): mark progress as "done"
$ws.Range("C139").Value = "done"

# Rows 180-190 (Binary Trees section): mark progress as "done"
$ws.Range("C180:C190").Value = "done"

# Reflect the author's last on-screen selection when the file was saved
$ws.Range("C140").Select()
